$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selection shown in the sheet view (B2:B32, active cell B2)
$ws.Range("B2:B32").Select()

# New duty-roster names for rows 2-32 (column B)
$names = @(
  "川田涼介",
  "豊島亮",
  "兒島大志郎",
  "日高泰聖",
  "Cox Matthew Jonah",
  "Hansen Jakob U",
  "Nicholas Tristan Aryasatyo",
  "小溝賢",
  "小野文哉",
  "渡部魁",
  "崎谷航平",
  "三神佳誠",
  "氏家琉貴",
  "羽賀尚生",
  "島田実",
  "足立耕平",
  "Yunjae",
  "神山修造",
  "志塚惇希",
  "川田涼介",
  "豊島亮",
  "兒島大志郎",
  "日高泰聖",
  "Cox Matthew Jonah",
  "Hansen Jakob U",
  "石井海成",
  "Nicholas Tristan Aryasatyo",
  "小溝賢",
  "小野文哉",
  "渡部魁",
  "崎谷航平"
)

for ($i = 0; $i -lt $names.Length; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 2).Value = $names[$i]
}
